$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Metadata sheet ("Metadata")
# ---------------------------------------------------------------

# Bump the Version and Date property values.
$ws1.Range("B3").Value = "2.0.2"
$ws1.Range("B8").Value = "2025-02-05T10:42:38+00:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws1.Range("A11:B11").EntireRow.Insert()

# The freshly inserted row doesn't inherit the surrounding data-row style,
# so copy the formatting (s="2") down from the row above it.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"
# Leave B11 blank (no Value assignment) - Jurisdiction has no value.

# ---------------------------------------------------------------
# Include sheet (was "Include from IPS.HAJJ.CONSENT")
# ---------------------------------------------------------------

$ws2.Name = "Include #0"

# Replace the 4 "CONSENT-xxx" concept rows (rows 2-5) with a single
# "All codes" row - drop rows 3,4,5 and repurpose row 2.
$ws2.Range("A3:B5").EntireRow.Delete()

# Row 1 becomes the "Codes" section header (single cell, A only).
$ws2.Range("A1").Value = "Codes"
$ws2.Range("B1").Clear()

# Row 2 becomes the "All codes" entry (single cell, A only).
$ws2.Range("A2").Value = "All codes"
$ws2.Range("B2").Clear()

# Rows 3 (blank separator) and 4 (System URI) keep their original content.
